$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("H2").Value = 2.9
$ws.Range("I2").Value = 2.63
$ws.Range("L2").Value = 3.6
$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5
$ws.Range("S2").Value = 2.88
$ws.Range("T2").Value = 1.4
$ws.Range("Y2").Value = 1.67
$ws.Range("Z2").Value = 2.1
$ws.Range("AD2").Value = 12
$ws.Range("AI2").Value = 5.5
$ws.Range("AO2").Value = 12
$ws.Range("AQ2").Value = 29

# Row 8 updates
$ws.Range("O8").Value = 1.29
$ws.Range("P8").Value = 3.5
$ws.Range("S8").Value = 1.88
$ws.Range("T8").Value = 1.93
